$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap the match data (columns F:V) between rows 55 and 56 ---
# Row 55 currently holds the "Pyunik Yerevan vs Ararat-Armenia" match,
# row 56 currently holds the "Van vs Alashkert" match; the edit swaps
# which row each match's data occupies (columns A:E - index/pais/torneio/
# temporada/data_partida - stay as-is since both rows share the same date).

$row55Vals = @(
    $ws.Cells.Item(55, 6).Value2,
    $ws.Cells.Item(55, 7).Value2,
    $ws.Cells.Item(55, 8).Value2,
    $ws.Cells.Item(55, 9).Value2,
    $ws.Cells.Item(55, 10).Value2,
    $ws.Cells.Item(55, 11).Value2,
    $ws.Cells.Item(55, 12).Value2,
    $ws.Cells.Item(55, 13).Value2,
    $ws.Cells.Item(55, 14).Value2,
    $ws.Cells.Item(55, 15).Value2,
    $ws.Cells.Item(55, 16).Value2,
    $ws.Cells.Item(55, 17).Value2,
    $ws.Cells.Item(55, 18).Value2,
    $ws.Cells.Item(55, 19).Value2,
    $ws.Cells.Item(55, 20).Value2,
    $ws.Cells.Item(55, 21).Value2,
    $ws.Cells.Item(55, 22).Value2
)

$row56Vals = @(
    $ws.Cells.Item(56, 6).Value2,
    $ws.Cells.Item(56, 7).Value2,
    $ws.Cells.Item(56, 8).Value2,
    $ws.Cells.Item(56, 9).Value2,
    $ws.Cells.Item(56, 10).Value2,
    $ws.Cells.Item(56, 11).Value2,
    $ws.Cells.Item(56, 12).Value2,
    $ws.Cells.Item(56, 13).Value2,
    $ws.Cells.Item(56, 14).Value2,
    $ws.Cells.Item(56, 15).Value2,
    $ws.Cells.Item(56, 16).Value2,
    $ws.Cells.Item(56, 17).Value2,
    $ws.Cells.Item(56, 18).Value2,
    $ws.Cells.Item(56, 19).Value2,
    $ws.Cells.Item(56, 20).Value2,
    $ws.Cells.Item(56, 21).Value2,
    $ws.Cells.Item(56, 22).Value2
)

for ($i = 0; $i -lt 17; $i++) {
    $ws.Cells.Item(55, 6 + $i).Value = $row56Vals[$i]
    $ws.Cells.Item(56, 6 + $i).Value = $row55Vals[$i]
}

# --- Step 2: append the new match row (row 71) ---
# Copy formatting from row 70 (the last existing row) so the new row's
# styles (A71 bold/centered/bordered index style, E71 date format) match.
$ws.Range("A70:V70").Copy()
$ws.Range("A71").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(71, 1).Value = 70
$ws.Cells.Item(71, 2).Value = "armenia"
$ws.Cells.Item(71, 3).Value = "premier-league"
$ws.Cells.Item(71, 4).Value = "2023-2024"
$ws.Cells.Item(71, 5).Value = 45233.625
$ws.Cells.Item(71, 6).Value = "Ararat Yerevan"
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = "Noah"
$ws.Cells.Item(71, 9).Value = 2
$ws.Cells.Item(71, 10).Value = 5.54
$ws.Cells.Item(71, 11).Value = "02/11/2023 03:12"
$ws.Cells.Item(71, 12).Value = 5.89
$ws.Cells.Item(71, 13).Value = "03/11/2023 14:59"
$ws.Cells.Item(71, 14).Value = 3.99
$ws.Cells.Item(71, 15).Value = "02/11/2023 03:12"
$ws.Cells.Item(71, 16).Value = 4.01
$ws.Cells.Item(71, 17).Value = "03/11/2023 14:59"
$ws.Cells.Item(71, 18).Value = 1.51
$ws.Cells.Item(71, 19).Value = "02/11/2023 03:12"
$ws.Cells.Item(71, 20).Value = 1.57
$ws.Cells.Item(71, 21).Value = "03/11/2023 14:56"
$ws.Cells.Item(71, 22).Value = "https://www.betexplorer.com/football/armenia/premier-league/ararat-yerevan-noah/23VPoXtP/"
